$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: recognizer name + URL, matching the existing rows' pattern.
$ws.Range("A6").Value = "auto_receipt"
$ws.Range("B6").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/general_receipt_recog"

# New hyperlink for the URL cell, mirroring the style applied to the other URL cells.
$ws.Hyperlinks.Add($ws.Range("B6"), "https://ocr-api.ccint.com/cci_ai/service/v1/general_receipt_recog")
$ws.Range("B6").Style = "Hyperlink"

# Column B widened to fit the new (longer) URL text, same "best fit" auto-sizing
# Excel applied to the existing rows.
$ws.Columns.Item(2).ColumnWidth = 52.166666666666664

# Match the diff's final selection state.
$ws.Range("C6:D6").Select()
